$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5997.8335
$ws.Range("I18").Value = 997.25
$ws.Range("J18").Value = 15999
$ws.Range("K18").Value = 997.25
$ws.Range("L18").Value = 15999
$ws.Range("M18").Value = -713.25
$ws.Range("N18").Value = -16567
$ws.Range("H39").Value = 1134.1666
$ws.Range("I39").Value = 561.3333
$ws.Range("K39").Value = 1683.9999
$ws.Range("M39").Value = -1387.9999
$ws.Range("H40").Value = 1910.0834
$ws.Range("I40").Value = 1849.091
$ws.Range("J40").Value = 1961.6923
$ws.Range("K40").Value = 1849.091
$ws.Range("L40").Value = 1961.6923
$ws.Range("M40").Value = -1674.091
$ws.Range("N40").Value = -2311.6923
$ws.Range("H116").Value = 461470.44
$ws.Range("I116").Value = 116556.2
$ws.Range("J116").Value = 892613.25
$ws.Range("K116").Value = 116556.2
$ws.Range("L116").Value = 892613.25
$ws.Range("M116").Value = -113114.2
$ws.Range("N116").Value = -899497.25
$ws.Range("H132").Value = 60378.523
$ws.Range("I132").Value = 68175.516
$ws.Range("K132").Value = 204526.548
$ws.Range("M132").Value = -201996.548
$ws.Range("H137").Value = 16951586
$ws.Range("J137").Value = 2944.5518
$ws.Range("L137").Value = 8833.6554
$ws.Range("N137").Value = -13933.6554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3336621.5
$ws.Range("I61").Value = 4765881.5
$ws.Range("K61").Value = 4765881.5
$ws.Range("M61").Value = -4765669.5
$ws.Range("H88").Value = 2550.111
$ws.Range("I88").Value = 1337.6666
$ws.Range("J88").Value = 3156.3333
$ws.Range("K88").Value = 1337.6666
$ws.Range("L88").Value = 3156.3333
$ws.Range("M88").Value = -931.6666
$ws.Range("N88").Value = -3968.3333
$ws.Range("H91").Value = 2550.111
$ws.Range("I91").Value = 1337.6666
$ws.Range("J91").Value = 3156.3333
$ws.Range("K91").Value = 1337.6666
$ws.Range("L91").Value = 3156.3333
$ws.Range("M91").Value = 66.33339999999998
$ws.Range("N91").Value = -5964.3333
$ws.Range("H102").Value = 72708.8
$ws.Range("I102").Value = 90261.25
$ws.Range("J102").Value = 2499
$ws.Range("K102").Value = 90261.25
$ws.Range("L102").Value = 2499
$ws.Range("M102").Value = -88639.25
$ws.Range("N102").Value = -5743
$ws.Range("H122").Value = 3127.5
$ws.Range("I122").Value = 2955.2856
$ws.Range("K122").Value = 8865.856800000001
$ws.Range("M122").Value = -6415.856800000001
$ws.Range("H132").Value = 557320.5600000001
$ws.Range("I132").Value = 668494.5600000001
$ws.Range("J132").Value = 1450.5
$ws.Range("K132").Value = 2005483.68
$ws.Range("L132").Value = 4351.5
$ws.Range("M132").Value = -2002953.68
$ws.Range("N132").Value = -9411.5
$ws.Range("H136").Value = 3336621.5
$ws.Range("I136").Value = 4765881.5
$ws.Range("K136").Value = 14297644.5
$ws.Range("M136").Value = -14295094.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2644.3635
$ws.Range("I105").Value = 2724.375
$ws.Range("J105").Value = 2431
$ws.Range("K105").Value = 2724.375
$ws.Range("L105").Value = 2431
$ws.Range("M105").Value = -977.375
$ws.Range("N105").Value = -5925
$ws.Range("H107").Value = 2373.795
$ws.Range("J107").Value = 4170.8184
$ws.Range("L107").Value = 4170.8184
$ws.Range("N107").Value = -8010.8184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17586.725
$ws.Range("I31").Value = 6560.722
$ws.Range("J31").Value = 35629.273
$ws.Range("K31").Value = 6560.722
$ws.Range("L31").Value = 35629.273
$ws.Range("M31").Value = -6265.722
$ws.Range("N31").Value = -36219.273
$ws.Range("H34").Value = 17586.725
$ws.Range("I34").Value = 6560.722
$ws.Range("J34").Value = 35629.273
$ws.Range("K34").Value = 6560.722
$ws.Range("L34").Value = 35629.273
$ws.Range("M34").Value = -6358.722
$ws.Range("N34").Value = -36033.273
$ws.Range("H99").Value = 2882.2727
$ws.Range("I99").Value = 2853.6667
$ws.Range("J99").Value = 2916.6
$ws.Range("K99").Value = 2853.6667
$ws.Range("L99").Value = 2916.6
$ws.Range("M99").Value = -1355.6667
$ws.Range("N99").Value = -5912.6
$ws.Range("H105").Value = 46806.75
$ws.Range("I105").Value = 52636.285
$ws.Range("K105").Value = 52636.285
$ws.Range("M105").Value = -50889.285
$ws.Range("H126").Value = 2882.2727
$ws.Range("I126").Value = 2853.6667
$ws.Range("J126").Value = 2916.6
$ws.Range("K126").Value = 8561.000100000001
$ws.Range("L126").Value = 8749.799999999999
$ws.Range("M126").Value = -6091.000100000001
$ws.Range("N126").Value = -13689.8
$ws.Range("H132").Value = 41877820
$ws.Range("I132").Value = 50001890
$ws.Range("K132").Value = 150005670
$ws.Range("M132").Value = -150003140
$ws.Range("H134").Value = 2345.037
$ws.Range("I134").Value = 2384.5
$ws.Range("J134").Value = 2171.4
$ws.Range("K134").Value = 7153.5
$ws.Range("L134").Value = 6514.200000000001
$ws.Range("M134").Value = -4618.5
$ws.Range("N134").Value = -11584.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5126.5
$ws.Range("J68").Value = 6495.1055
$ws.Range("L68").Value = 19485.3165
$ws.Range("N68").Value = -21107.3165
$ws.Range("H71").Value = 5126.5
$ws.Range("J71").Value = 6495.1055
$ws.Range("L71").Value = 58455.9495
$ws.Range("N71").Value = -66567.94949999999
$ws.Range("H86").Value = 751
$ws.Range("I86").Value = 751
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2253
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1067
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 751
$ws.Range("I89").Value = 751
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6759
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -831
$ws.Range("N89").ClearContents()
$ws.Range("H137").Value = 3021.7
$ws.Range("I137").Value = 1026.8572
$ws.Range("J137").Value = 7676.3335
$ws.Range("K137").Value = 3080.5716
$ws.Range("L137").Value = 23029.0005
$ws.Range("M137").Value = 2019.4284
$ws.Range("N137").Value = -33229.00049999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4527.952
$ws.Range("I122").Value = 3685.7334
$ws.Range("K122").Value = 11057.2002
$ws.Range("M122").Value = -8607.200199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1659.3
$ws.Range("I61").Value = 1524.625
$ws.Range("J61").Value = 2198
$ws.Range("K61").Value = 1524.625
$ws.Range("L61").Value = 2198
$ws.Range("M61").Value = -1322.625
$ws.Range("N61").Value = -2602
$ws.Range("H113").Value = 1659.3
$ws.Range("I113").Value = 1524.625
$ws.Range("J113").Value = 2198
$ws.Range("K113").Value = 1524.625
$ws.Range("L113").Value = 2198
$ws.Range("M113").Value = 645.375
$ws.Range("N113").Value = -6538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11930557
$ws.Range("I136").Value = 13361424
$ws.Range("J136").Value = 6666.3335
$ws.Range("K136").Value = 40084272
$ws.Range("L136").Value = 19999.0005
$ws.Range("M136").Value = -25099.0005
